# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colour scheme (used by the notes master)
#   ppt/theme/theme2.xml -> "Integral" colour scheme (used by the slide master / every slide)
#
# The authored edit swaps the two themes' contents: the slide master (and
# therefore every slide) switches from the "Integral" palette to the
# "Office Theme" palette, while the notes master ends up with the
# "Integral" palette. The only real content difference between the two
# theme parts in this deck is the 12-colour scheme (fonts/format scheme are
# identical), so we reproduce the swap by rewriting the live theme's colour
# scheme - i.e. the one bound to $p.SlideMaster - to the "Office Theme"
# values, the same way a user would via Design > Variants > Colors in the
# PowerPoint UI (COM: ThemeColorScheme.Colors(index).RGB).

function ToRgbLong($hex) {
    $r = ($hex -band 0xFF0000) -shr 16
    $g = ($hex -band 0x00FF00) -shr 8
    $b = ($hex -band 0x0000FF)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Index order exposed by ThemeColorScheme.Colors(): 1-12 =>
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $colorScheme.Colors($i).RGB = ToRgbLong $officeThemeColors[$i - 1]
}
